$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "AKBANK" column (D) benchmark figures for rows 3-14
# (D7 and D12-adjacent cells not in this set are intentionally left alone)
$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("D11").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("D14").ClearContents()

# Clear the "DENIZBANK" column (J) benchmark figures for rows 3-14
$ws.Range("J3").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("J6").ClearContents()
$ws.Range("J8").ClearContents()
$ws.Range("J9").ClearContents()
$ws.Range("J10").ClearContents()
$ws.Range("J11").ClearContents()
$ws.Range("J13").ClearContents()
$ws.Range("J14").ClearContents()

# Updated figures for GELEN SWIFT row (13): İŞBANKASI (E) and FINASNBANK (K)
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 851,5 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"
